# agora +1 forma de pagamento pede o valor exato
# Add a new worksheet "Contas Fechamento" right after the existing "Sheet",
# with the same header row (DATA / CONTAS / VALOR) used on the first sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new worksheet immediately after the first sheet.
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Contas Fechamento"

# Header row, matching the layout of the original sheet.
$newSheet.Range("A1").Value = "DATA"
$newSheet.Range("B1").Value = "CONTAS"
$newSheet.Range("C1").Value = "VALOR"
